# Updated metadata values produced by the improved corner-detection /
# dominant-color algorithm. Row 2 (frame 0) no longer has meaningful
# detections so it is zeroed out; rows 3 and 4 (frames 1 and 2) pick up
# refreshed measurements.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (frame 0) - no corners/ellipse detected anymore, clear the metrics
$ws.Range("B2:X2").Value = 0
$ws.Range("Z2:AA2").Value = 0

# Row 3 (frame 1) - refreshed dominant color / corner metrics
$ws.Range("B3").Value = 73.62410059253746
$ws.Range("C3").Value = 181.1961237301674
$ws.Range("D3").Value = 146.3996190745278
$ws.Range("H3").Value = 86.60545657905126
$ws.Range("I3").Value = 38.43600703045831
$ws.Range("J3").Value = 123.4291763202008
$ws.Range("N3").Value = 38.27906822842785
$ws.Range("O3").Value = 14.65614453049086
$ws.Range("P3").Value = 93.62921050889517

# Row 4 (frame 2) - refreshed dominant color / corner metrics
$ws.Range("B4").Value = 109.6901181963495
$ws.Range("C4").Value = 97.67287850589261
$ws.Range("D4").Value = 186.4254334410942
$ws.Range("E4").Value = 87.42396519736401
$ws.Range("F4").Value = 104.9072995969533
$ws.Range("G4").Value = 154.5306762203259
$ws.Range("H4").Value = 33.1519060454403
$ws.Range("I4").Value = 43.12033115132662
$ws.Range("J4").Value = 82.43357720447598
$ws.Range("K4").Value = 74.51240661686211
$ws.Range("L4").Value = 81.89701173959355
$ws.Range("M4").Value = 127.7704108858022
$ws.Range("N4").Value = 50.21404158360141
$ws.Range("O4").Value = 59.29592708629858
$ws.Range("P4").Value = 109.0707775562541
